$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'69.818.87"
$ws.Range("D3").Formula = "'3.527.62"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Formula = "'605.14"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Formula = "'195.26"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Formula = "'0.202"
$ws.Range("E9").Value = "  -5.12%  "
$ws.Range("D10").Formula = "'0.645"
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("D11").Formula = "'53.37"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Formula = "'0.0000302"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Formula = "'9.46"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Formula = "'4.093.72"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Formula = "'594.01"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Formula = "'69.932.10"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Formula = "'12.70"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Formula = "'18.96"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Formula = "'0.122"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Formula = "'3.514.03"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Formula = "'0.983"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Formula = "'17.80"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Formula = "'103.36"
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").Formula = "'5.15"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Formula = "'3.04"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Formula = "'10.77"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Formula = "'9.51"
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("D29").Formula = "'33.20"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("D30").Formula = "'7.04"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Formula = "'4.21"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Formula = "'12.31"
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Formula = "'63.45"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Formula = "'3.17"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Formula = "'3.784.39"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Formula = "'1.00"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Formula = "'0.0₃0809"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").Formula = "'510.78"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("D40").Formula = "'0.390"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Formula = "'36.40"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("D44").Formula = "'0.0448"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").Formula = "'2.81"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("D49").Formula = "'8.46"
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Formula = "'1.34"
$ws.Range("E50").Value = "  +3.31%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Formula = "'0.000248"
$ws.Range("E51").Value = "  +4.79%  "
